$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the record for RICARDO RAFAEL HERRERA DE LAVALLE (row 16).
# This shifts every row below it up by one, which naturally moves the
# "last row" special border styling from the old row 22 to the new
# last data row (21), and the signature block from rows 27/28 to 26/27.
$ws.Rows("16:16").Delete()

# After the shift, the four ROMEL ALEXANDER SANCHEZ LONDOÃ?O period rows
# (now rows 18-21) carry their old descending period order (2506..2503).
# Put the periods (and matching Valor Mora) back into ascending order
# 2503..2506, updating both the period label (col E) and Valor Mora (col F).
$ws.Range("E18").Value = "2503"
$ws.Range("F18").Value = 34635

$ws.Range("E19").Value = "2504"
$ws.Range("F19").Value = 64940

$ws.Range("E20").Value = "2505"
$ws.Range("F20").Value = 64940

$ws.Range("E21").Value = "2506"
$ws.Range("F21").Value = 64940

# Update the summary figures to reflect the new data set.
$ws.Range("E11").Value = 336079
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 5
